$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores locale-formatted numeric-looking text
# (e.g. "63.234.50", using "." as both thousands separator and
# decimal point, or long decimals like "0.0000236"). Force each
# updated Price cell to Text format before writing its new value so
# Excel keeps the literal digits instead of reinterpreting the text
# as a floating point number (which would round "588.29" to
# "588.28999999999996", or corrupt multi-dot values like
# "63.335.22").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.335.22'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.056.78'
$ws.Range("E3").Value = '  -2.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.29'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.38'
$ws.Range("E6").Value = '  +6.08%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.054.67'
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.82'
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.03'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000237'
$ws.Range("E14").Value = '  -3.89%  '
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.559.03'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.334.61'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.12'
$ws.Range("E18").Value = '  -2.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.053.21'
$ws.Range("E19").Value = '  -2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '472.53'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.34'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").Value = '  -3.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.50'
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.59'
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.79'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.39'
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("E29").Value = '  +2.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -2.85%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("E34").Value = '  -2.33%  '
$ws.Range("E35").Value = '  -3.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("E37").Value = '  +0.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.97'
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("E39").Value = '  -3.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.25'
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.72'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '441.20'
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.288'
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.56'
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0358'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.787.66'
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.96'
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.05'
$ws.Range("E50").Value = '  +4.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.23'
$ws.Range("E51").Value = '  +0.26%  '
